$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old hyperlinks (previously anchored on column D: D4, D5, D6) ---
$ws.Range("D4:D6").Hyperlinks.Delete()

# --- The old hyperlink cells carried the "Hyperlink" cell style; reset it back to
#     Normal since column D no longer holds a hyperlink (it becomes "provider"). ---
$ws.Range("D4:D6").Style = "Normal"

# --- Clear the old cell contents (keeps the existing header row styling intact) ---
$ws.UsedRange.ClearContents()

# --- Header row (column order: company_name, career_url, keywords, provider, slug) ---
$ws.Range("A1").Value = "company_name"
$ws.Range("B1").Value = "career_url"
$ws.Range("C1").Value = "keywords"
$ws.Range("D1").Value = "provider"
$ws.Range("E1").Value = "slug"

# --- Data rows ---
# Row 2: OpenAI
$ws.Range("A2").Value = "OpenAI"
$ws.Range("B2").Value = "https://boards.greenhouse.io/openai"
$ws.Range("C2").Value = "LLM; research; machine learning; NLP"
$ws.Range("D2").Value = "greenhouse"
$ws.Range("E2").Value = "openai"
$ws.Hyperlinks.Add($ws.Range("B2"), "https://boards.greenhouse.io/openai") | Out-Null

# Row 3: Databricks
$ws.Range("A3").Value = "Databricks"
$ws.Range("B3").Value = "https://jobs.lever.co/databricks"
$ws.Range("C3").Value = "LLM; NLP; data scientist; ML engineer"
$ws.Range("D3").Value = "lever"
$ws.Range("E3").Value = "databricks"
$ws.Hyperlinks.Add($ws.Range("B3"), "https://jobs.lever.co/databricks") | Out-Null

# Row 4: Canva
$ws.Range("A4").Value = "Canva"
$ws.Range("B4").Value = "https://www.smartrecruiters.com/Canva"
$ws.Range("C4").Value = "machine learning; AI engineer; research"
$ws.Range("D4").Value = "smartrecruiters"
$ws.Range("E4").Value = "Canva"
$ws.Hyperlinks.Add($ws.Range("B4"), "https://www.smartrecruiters.com/Canva") | Out-Null

# Row 5: Stripe
$ws.Range("A5").Value = "Stripe"
$ws.Range("B5").Value = "https://boards.greenhouse.io/stripe"
$ws.Range("C5").Value = "LLM; data scientist; machine learning"
$ws.Range("D5").Value = "greenhouse"
$ws.Range("E5").Value = "stripe"
$ws.Hyperlinks.Add($ws.Range("B5"), "https://boards.greenhouse.io/stripe") | Out-Null

# Row 6: NVIDIA
$ws.Range("A6").Value = "NVIDIA"
$ws.Range("B6").Value = "https://nvidia.wd5.myworkdayjobs.com/en-US/NVIDIAExternalCareerSite"
$ws.Range("C6").Value = "NLP; LLM; machine learning; research"
$ws.Range("D6").Value = "workday"
$ws.Hyperlinks.Add($ws.Range("B6"), "https://nvidia.wd5.myworkdayjobs.com/en-US/NVIDIAExternalCareerSite") | Out-Null
